# Update IESO report from GitHub Actions
# Refresh the report timestamp and the updated LMP / price figures produced
# by the latest data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-05-04T03:09:54"
$ws.Range("G4").Value = 167.31
$ws.Range("H4").Value = 164
$ws.Range("S4").Value = 7.65
$ws.Range("T4").Value = 7.56
$ws.Range("U4").Value = 7.55
$ws.Range("V4").Value = 7.7
$ws.Range("S6").Value = 0.21
$ws.Range("T6").Value = 0.2
$ws.Range("U6").Value = 0.2
$ws.Range("V6").Value = 0.18
$ws.Range("G7").Value = 167.31
$ws.Range("H7").Value = 164
$ws.Range("H9").Value = 16.76
$ws.Range("S9").Value = 7.98
$ws.Range("T9").Value = 7.68
$ws.Range("U9").Value = 7.66
$ws.Range("V9").Value = 7.83
$ws.Range("S11").Value = 0.54
$ws.Range("T11").Value = 0.32
$ws.Range("V11").Value = 0.31
$ws.Range("H12").Value = 16.76
$ws.Range("H14").Value = 22.6
$ws.Range("S14").Value = 7.98
$ws.Range("T14").Value = 7.68
$ws.Range("U14").Value = 7.66
$ws.Range("V14").Value = 7.83
$ws.Range("S16").Value = 0.54
$ws.Range("T16").Value = 0.32
$ws.Range("V16").Value = 0.31
$ws.Range("H17").Value = 16.76
$ws.Range("H19").Value = 38.5
$ws.Range("S19").Value = 7.72
$ws.Range("T19").Value = 7.62
$ws.Range("U19").Value = 114.94
$ws.Range("V19").Value = 114.94
$ws.Range("W19").Value = 42.95
$ws.Range("S21").Value = 0.29
$ws.Range("T21").Value = 0.27
$ws.Range("U21").Value = 0.25
$ws.Range("H22").Value = 38.5
$ws.Range("U22").Value = 107.34
$ws.Range("V22").Value = 107.19
$ws.Range("G24").Value = 167.31
$ws.Range("H24").Value = 164
$ws.Range("S24").Value = 7.72
$ws.Range("T24").Value = 7.62
$ws.Range("U24").Value = 7.61
$ws.Range("V24").Value = 7.75
$ws.Range("W24").Value = 6.54
$ws.Range("S26").Value = 0.29
$ws.Range("T26").Value = 0.27
$ws.Range("U26").Value = 0.25
$ws.Range("G27").Value = 167.31
$ws.Range("H27").Value = 164
$ws.Range("H29").Value = 38.5
$ws.Range("S29").Value = 7.78
$ws.Range("T29").Value = 7.67
$ws.Range("S31").Value = 0.34
$ws.Range("T31").Value = 0.31
$ws.Range("U31").Value = 0.31
$ws.Range("V31").Value = 0.3
$ws.Range("H32").Value = 38.5
$ws.Range("U32").Value = 107.34
$ws.Range("V32").Value = 107.19
$ws.Range("H34").Value = 25
$ws.Range("M34").Value = 6.03
$ws.Range("N34").Value = 6.12
$ws.Range("S34").Value = 8.08
$ws.Range("T34").Value = 7.65
$ws.Range("U34").Value = 7.64
$ws.Range("V34").Value = 7.82
$ws.Range("W34").Value = 6.83
$ws.Range("S36").Value = 0.64
$ws.Range("T36").Value = 0.3
$ws.Range("U36").Value = 0.28
$ws.Range("V36").Value = 0.3
$ws.Range("H37").Value = 19.16
$ws.Range("M37").Value = 0.62
$ws.Range("N37").Value = 0.71
$ws.Range("G39").Value = 167.31
$ws.Range("H39").Value = 164
$ws.Range("S39").Value = 7.65
$ws.Range("T39").Value = 7.56
$ws.Range("U39").Value = 7.55
$ws.Range("V39").Value = 7.7
$ws.Range("S41").Value = 0.21
$ws.Range("T41").Value = 0.2
$ws.Range("U41").Value = 0.2
$ws.Range("V41").Value = 0.18
$ws.Range("G42").Value = 167.31
$ws.Range("H42").Value = 164
$ws.Range("S44").Value = 7.42
$ws.Range("T44").Value = 7.36
$ws.Range("U44").Value = 7.36
$ws.Range("V44").Value = 7.51
$ws.Range("S46").Value = -0.02
$ws.Range("V46").Value = 0
$ws.Range("S49").Value = 6.78
$ws.Range("T49").Value = 6.72
$ws.Range("U49").Value = 6.67
$ws.Range("V49").Value = 6.86
$ws.Range("W49").Value = 5.73
$ws.Range("S51").Value = -0.66
$ws.Range("T51").Value = -0.63
$ws.Range("U51").Value = -0.6899999999999999
$ws.Range("V51").Value = -0.65
$ws.Range("S54").Value = 6.73
$ws.Range("T54").Value = 6.65
$ws.Range("U54").Value = 6.77
$ws.Range("V54").Value = 6.92
$ws.Range("S56").Value = -0.71
$ws.Range("T56").Value = -0.7
$ws.Range("U56").Value = -0.59
$ws.Range("V56").Value = -0.59
$ws.Range("W56").Value = -0.5
$ws.Range("S59").Value = 7.57
$ws.Range("T59").Value = 7.51
$ws.Range("U59").Value = 7.5
$ws.Range("V59").Value = 7.66
$ws.Range("W59").Value = 6.42
$ws.Range("S61").Value = 0.13
$ws.Range("T61").Value = 0.15
$ws.Range("U61").Value = 0.15
$ws.Range("V61").Value = 0.15
$ws.Range("S64").Value = 7.66
$ws.Range("T64").Value = 7.6
$ws.Range("U64").Value = 7.6
$ws.Range("V64").Value = 7.76
$ws.Range("W64").Value = 6.5
$ws.Range("S66").Value = 0.22
$ws.Range("T66").Value = 0.24
$ws.Range("U66").Value = 0.24
$ws.Range("V66").Value = 0.25
$ws.Range("W66").Value = 0.19
$ws.Range("H69").Value = 15.02
$ws.Range("S69").Value = 7.61
$ws.Range("T69").Value = 7.54
$ws.Range("U69").Value = 7.54
$ws.Range("V69").Value = 7.7
$ws.Range("W69").Value = 6.45
$ws.Range("S71").Value = 0.18
$ws.Range("T71").Value = 0.19
$ws.Range("U71").Value = 0.19
$ws.Range("H72").Value = 15.02
$ws.Range("S74").Value = 7.1
$ws.Range("T74").Value = 7.02
$ws.Range("U74").Value = 7.02
$ws.Range("V74").Value = 7.18
$ws.Range("S76").Value = -0.33
$ws.Range("T76").Value = -0.33
$ws.Range("U76").Value = -0.33
$ws.Range("V76").Value = -0.33
$ws.Range("S79").Value = 7.44
$ws.Range("T79").Value = 7.36
$ws.Range("U79").Value = 7.35
$ws.Range("V79").Value = 7.51
$ws.Range("W79").Value = 6.3
$ws.Range("S84").Value = 6.66
$ws.Range("T84").Value = 6.58
$ws.Range("U84").Value = 6.89
$ws.Range("V84").Value = 7.07
$ws.Range("S86").Value = -0.78
$ws.Range("T86").Value = -0.77
$ws.Range("U86").Value = -0.46
$ws.Range("V86").Value = -0.45
$ws.Range("H89").Value = 38.5
$ws.Range("S89").Value = 7.78
$ws.Range("T89").Value = 7.67
$ws.Range("S91").Value = 0.34
$ws.Range("T91").Value = 0.31
$ws.Range("U91").Value = 0.31
$ws.Range("V91").Value = 0.3
$ws.Range("H92").Value = 38.5
$ws.Range("U92").Value = 107.34
$ws.Range("V92").Value = 107.19
